$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the extra "or if u round the .66 ..." run (and with it,
#    the _GoBack bookmark that used to sit right after it) from the
#    paragraph that ends in ".4995".
# ------------------------------------------------------------------
$d.Content.Find.Execute(" or if u round the .66 to .67 and u get .5000025", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 2) Add a brand-new paragraph after the "Bonus: ..." paragraph with
#    two runs:
#      - "See above please" (yellow highlight)
#      - " we did this up there" (no highlight)
#    and put the _GoBack bookmark at the very end of it.
# ------------------------------------------------------------------
$bonusText = "Bonus: Can you find a split choice where GINI and entropy select different splits?"
$full = $d.Content.Text
$bonusIdx = $full.IndexOf($bonusText)
$bonusEnd = $bonusIdx + $bonusText.Length

$insertPoint = $d.Range($bonusEnd, $bonusEnd)
$insertPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$pr = $newPara.Range

# Insert both runs of text in one go (using a placeholder token for the
# first run), plus a temporary trailing token. Keeping a little bit of
# extra text after the real content means later range math never has
# to land exactly on the last character of the document, a position
# this host mis-resolves (a collapsed Range there snaps back to the
# start of the document instead of staying put).
$pr.InsertAfter("PLACEHOLDERONE we did this up there PLACEHOLDERTWO")

# Apply the yellow highlight to just the first run via Find/Replace
# formatting - this correctly splits the run instead of touching the
# whole paragraph the way directly setting HighlightColorIndex on a
# sub-range does.
$find = $pr.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Highlight = $true
$find.Execute("PLACEHOLDERONE", $false, $false, $false, $false, $false, $true, 1, $false, "See above please", 2)

# Place the _GoBack bookmark right before the trailing placeholder
# (i.e. exactly at the end of the real text), then delete the
# placeholder text (and the space in front of it) so only the real
# wording remains.
$full2 = $d.Content.Text
$markPos = $full2.IndexOf(" PLACEHOLDERTWO")
$bmRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$full3 = $d.Content.Text
$markPos2 = $full3.IndexOf(" PLACEHOLDERTWO")
$markRange = $d.Range($markPos2, $markPos2 + " PLACEHOLDERTWO".Length)
$markRange.Delete()
